# Auto-generated Excel COM-interop script applying the diff to Leviathan_Profits workbook.
# For each affected sheet, updates cell values (and clears cells that were removed).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 46022.332
$ws.Range("I47").Value = 46022.332
$ws.Range("K47").Value = 46022.332
$ws.Range("M47").Value = -45050.332
$ws.Range("H88").Value = 202250.8
$ws.Range("I88").Value = 2700
$ws.Range("K88").Value = 2700
$ws.Range("M88").Value = -2294
$ws.Range("H91").Value = 202250.8
$ws.Range("I91").Value = 2700
$ws.Range("K91").Value = 2700
$ws.Range("M91").Value = -1296
$ws.Range("H135").Value = 655.04
$ws.Range("I135").Value = 579.9524
$ws.Range("K135").Value = 5219.5716
$ws.Range("M135").Value = -2684.5716
$ws.Range("H137").Value = 127162.125
$ws.Range("I137").Value = 2799.5
$ws.Range("J137").Value = 251524.75
$ws.Range("K137").Value = 8398.5
$ws.Range("L137").Value = 754574.25
$ws.Range("M137").Value = -5848.5
$ws.Range("N137").Value = -759674.25
$ws.Range("H138").Value = 2076.4443
$ws.Range("I138").Value = 1337.8667
$ws.Range("K138").Value = 4013.6001
$ws.Range("M138").Value = 1126.3999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 128705.58
$ws.Range("I32").Value = 134705.9
$ws.Range("J32").Value = 101704.1
$ws.Range("K32").Value = 134705.9
$ws.Range("L32").Value = 101704.1
$ws.Range("M32").Value = -134418.9
$ws.Range("N32").Value = -102278.1
$ws.Range("H45").Value = 16009.435
$ws.Range("I45").Value = 16689.072
$ws.Range("K45").Value = 16689.072
$ws.Range("M45").Value = -16312.072
$ws.Range("H48").Value = 125099.5
$ws.Range("J48").Value = 125099.5
$ws.Range("L48").Value = 125099.5
$ws.Range("N48").Value = -125867.5
$ws.Range("H61").Value = 3115.9473
$ws.Range("I61").Value = 3066.8333
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 3066.8333
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2854.8333
$ws.Range("N61").Value = -4424
$ws.Range("H80").Value = 19998.334
$ws.Range("J80").Value = 19998.334
$ws.Range("L80").Value = 19998.334
$ws.Range("N80").Value = -21994.334
$ws.Range("H83").Value = 19998.334
$ws.Range("J83").Value = 19998.334
$ws.Range("L83").Value = 59995.00199999999
$ws.Range("N83").Value = -69979.00199999999
$ws.Range("H136").Value = 3115.9473
$ws.Range("I136").Value = 3066.8333
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 9200.499899999999
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -6650.499899999999
$ws.Range("N136").Value = -17100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 872.6667
$ws.Range("J80").Value = 1060.5
$ws.Range("L80").Value = 1060.5
$ws.Range("N80").Value = -3056.5
$ws.Range("H83").Value = 872.6667
$ws.Range("J83").Value = 1060.5
$ws.Range("L83").Value = 5302.5
$ws.Range("N83").Value = -15286.5
$ws.Range("H107").Value = 7284.143
$ws.Range("I107").Value = 7998.1665
$ws.Range("K107").Value = 7998.1665
$ws.Range("M107").Value = -6078.1665
$ws.Range("H134").Value = 2717.111
$ws.Range("I134").Value = 2805
$ws.Range("J134").Value = 2014
$ws.Range("K134").Value = 8415
$ws.Range("L134").Value = 6042
$ws.Range("M134").Value = -5880
$ws.Range("N134").Value = -11112

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2029.091
$ws.Range("I31").Value = 1630.3889
$ws.Range("K31").Value = 1630.3889
$ws.Range("M31").Value = -1335.3889
$ws.Range("H34").Value = 2029.091
$ws.Range("I34").Value = 1630.3889
$ws.Range("K34").Value = 1630.3889
$ws.Range("M34").Value = -1428.3889
$ws.Range("H58").Value = 6984.2104
$ws.Range("I58").Value = 7356.4375
$ws.Range("J58").Value = 4999
$ws.Range("K58").Value = 7356.4375
$ws.Range("L58").Value = 4999
$ws.Range("M58").Value = -7153.4375
$ws.Range("N58").Value = -5405
$ws.Range("H99").Value = 2900
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H126").Value = 2900
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 2300.5
$ws.Range("I132").Value = 2099
$ws.Range("K132").Value = 6297
$ws.Range("M132").Value = -3767
$ws.Range("H134").Value = 2549.6
$ws.Range("I134").Value = 2361.7097
$ws.Range("K134").Value = 7085.1291
$ws.Range("M134").Value = -4550.1291
$ws.Range("H136").Value = 6984.2104
$ws.Range("I136").Value = 7356.4375
$ws.Range("J136").Value = 4999
$ws.Range("K136").Value = 22069.3125
$ws.Range("L136").Value = 14997
$ws.Range("M136").Value = -19519.3125
$ws.Range("N136").Value = -20097

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 5147.4
$ws.Range("I46").Value = 295.8
$ws.Range("J46").Value = 9999
$ws.Range("K46").Value = 887.4000000000001
$ws.Range("L46").Value = 29997
$ws.Range("M46").Value = -796.4000000000001
$ws.Range("N46").Value = -30179

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 80000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 80000
$ws.Range("K42").Value = 0
$ws.Range("L42").ClearContents()
$ws.Range("M42").Value = 80000
$ws.Range("N42").Value = -80970
$ws.Range("H46").Value = 14812.25
$ws.Range("J46").Value = 49999
$ws.Range("L46").Value = 49999
$ws.Range("N46").Value = -50311
$ws.Range("H115").Value = 80000
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 80000
$ws.Range("K115").Value = 0
$ws.Range("L115").ClearContents()
$ws.Range("M115").Value = 80000
$ws.Range("N115").Value = -82350
$ws.Range("H132").Value = 2372.25
$ws.Range("I132").Value = 2495
$ws.Range("J132").Value = 2249.5
$ws.Range("K132").Value = 7485
$ws.Range("L132").Value = 6748.5
$ws.Range("M132").Value = -4955
$ws.Range("N132").Value = -11808.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6665.7856
$ws.Range("I40").Value = 8834.5
$ws.Range("J40").Value = 5039.25
$ws.Range("K40").Value = 8834.5
$ws.Range("L40").Value = 5039.25
$ws.Range("M40").Value = -8698.5
$ws.Range("N40").Value = -5311.25
$ws.Range("H68").Value = 2837.923
$ws.Range("I68").Value = 2343.889
$ws.Range("J68").Value = 3949.5
$ws.Range("K68").Value = 2343.889
$ws.Range("L68").Value = 3949.5
$ws.Range("M68").Value = -1594.889
$ws.Range("N68").Value = -5447.5
$ws.Range("H71").Value = 2837.923
$ws.Range("I71").Value = 2343.889
$ws.Range("J71").Value = 3949.5
$ws.Range("K71").Value = 11719.445
$ws.Range("L71").Value = 19747.5
$ws.Range("M71").Value = -7975.445
$ws.Range("N71").Value = -27235.5
$ws.Range("H82").Value = 930.75
$ws.Range("I82").Value = 980.2
$ws.Range("K82").Value = 980.2
$ws.Range("M82").Value = -619.2
$ws.Range("H85").Value = 930.75
$ws.Range("I85").Value = 980.2
$ws.Range("K85").Value = 980.2
$ws.Range("M85").Value = 267.8
$ws.Range("H132").Value = 7278.2
$ws.Range("I132").Value = 7553.9443
$ws.Range("J132").Value = 4796.5
$ws.Range("K132").Value = 22661.8329
$ws.Range("L132").Value = 14389.5
$ws.Range("M132").Value = -20131.8329
$ws.Range("N132").Value = -19449.5
$ws.Range("H136").Value = 3048.1
$ws.Range("I136").Value = 2435.625
$ws.Range("K136").Value = 7306.875
$ws.Range("M136").Value = -4756.875

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 30389
$ws.Range("J93").Value = 30389
$ws.Range("L93").Value = 30389
$ws.Range("N93").Value = -35381
$ws.Range("H113").Value = 1700
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1700
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").Value = 5100
$ws.Range("N113").Value = -9440
$ws.Range("H132").Value = 6171.0884
$ws.Range("I132").Value = 8253.608
$ws.Range("J132").Value = 1816.7273
$ws.Range("K132").Value = 24760.824
$ws.Range("L132").Value = 5450.1819
$ws.Range("M132").Value = -22230.824
$ws.Range("N132").Value = -10510.1819
$ws.Range("H136").Value = 1789.1
$ws.Range("I136").Value = 1789.1
$ws.Range("K136").Value = 5367.299999999999
$ws.Range("M136").Value = -2817.299999999999

Write-Host "Applied all Leviathan_Profits cell updates."